$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) values for the refreshed crypto snapshot.
# Column D values are forced to Text format first so Excel stores them as strings
# (matching the workbook author's original inline-string cells) instead of parsing them as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.577.60"
$ws.Range("E2").Value = "  +1.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.03"
$ws.Range("E3").Value = "  +1.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.27%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.21"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5344"
$ws.Range("E7").Value = "  +0.74%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3979"
$ws.Range("E8").Value = "  +5.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07711"
$ws.Range("E9").Value = "  +3.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.119"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.86"
$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.346"
$ws.Range("E12").Value = "  +2.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.01"
$ws.Range("E13").Value = "  +2.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.600"
$ws.Range("E14").Value = "  +3.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.003"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.831.29"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.11"
$ws.Range("E17").Value = "  +2.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001088"
$ws.Range("E18").Value = "  +1.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06577"
$ws.Range("E19").Value = "  +1.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.77"
$ws.Range("E20").Value = "  +3.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.086"
$ws.Range("E22").Value = "  +2.75%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.604.60"
$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.20"
$ws.Range("E24").Value = "  +0.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.247"
$ws.Range("E25").Value = "  +7.85%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.75"
$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.041.95"
$ws.Range("E27").Value = "  +1.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.14"
$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.422"
$ws.Range("E29").Value = "  +3.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.99"
$ws.Range("E30").Value = "  +2.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.138"
$ws.Range("E31").Value = "  +1.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1122"
$ws.Range("E32").Value = "  +1.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.739"
$ws.Range("E33").Value = "  +2.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.663"
$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07288"
$ws.Range("E35").Value = "  +0.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2259"
$ws.Range("E36").Value = "  +1.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02348"
$ws.Range("E37").Value = "  +2.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.908"
$ws.Range("E38").Value = "  +4.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.206"
$ws.Range("E39").Value = "  +2.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.41"
$ws.Range("E40").Value = "  +2.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6300"
$ws.Range("E41").Value = "  +2.10%  "

$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.392"
$ws.Range("E44").Value = "  -2.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.59"
$ws.Range("E45").Value = "  +1.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5899"
$ws.Range("E46").Value = "  +2.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.721"
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.34"
$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.997"
$ws.Range("E49").Value = "  +3.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.197"
$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06937"
$ws.Range("E51").Value = "  +1.70%  "
